$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item("removed").Name = "removed row"
$wb.Worksheets.Item("added").Name = "added row"
$wb.Worksheets.Item("changed").Name = "changed row"

# Highlight the changed cells with an orange background
$ws = $wb.Worksheets.Item("changed row")
$changedCells = @("AN2", "AE3", "AN3", "AN4", "D5")
foreach ($cellRef in $changedCells) {
    $ws.Range($cellRef).Interior.Color = 0x00A5FF
}
